# Fix field validation and UI guide in Excel upload
# Update the "Input" test sheet to match the standard PO template:
# new column layout (발주일자, 납기일자, 거래처명, 거래처 이메일, 납품처명,
# 납품처 이메일, 프로젝트명, 대분류, 중분류, 소분류, 품목명, 규격, 수량,
# 단가, 총금액, 비고) and drop the old header bold/border styling.
# Also tidy the 갑지/을지 sheets by removing the stray empty 비고 cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Input": rebuild header + data with the new standard layout
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Input")

# Wipe out the old content (values + the bold/bordered header formatting)
$ws.Range("A1:Q3").Clear()

$headers = @(
    "발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명",
    "납품처 이메일", "프로젝트명", "대분류", "중분류", "소분류",
    "품목명", "규격", "수량", "단가", "총금액", "비고"
)
for ($col = 0; $col -lt $headers.Length; $col++) {
    $ws.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# Row 2 data
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-09-10"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2025-09-02"
$ws.Cells.Item(2, 3).Value = "티에스이앤씨"
$ws.Cells.Item(2, 4).Value = "티에스이앤씨@example.com"
$ws.Cells.Item(2, 5).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(2, 6).Value = "delivery@example.com"
$ws.Cells.Item(2, 7).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(2, 8).Value = "2. 부자재비"
$ws.Cells.Item(2, 9).Value = "2) 창호"
$ws.Cells.Item(2, 10).Value = "기타"
$ws.Cells.Item(2, 11).Value = "3차 - 스크류 (둥근머리 8*25)"
$ws.Cells.Item(2, 12).Value = "KS규격-1"
$ws.Cells.Item(2, 13).Value = 500
$ws.Cells.Item(2, 14).Value = 19
$ws.Cells.Item(2, 15).Value = 10450

# Row 3 data
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "2025-08-28"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2025-09-05"
$ws.Cells.Item(3, 3).Value = "티에스이앤씨"
$ws.Cells.Item(3, 4).Value = "티에스이앤씨@example.com"
$ws.Cells.Item(3, 5).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(3, 6).Value = "delivery@example.com"
$ws.Cells.Item(3, 7).Value = "힐스테이트 도곡동1차"
$ws.Cells.Item(3, 8).Value = "6. 안전관리비"
$ws.Cells.Item(3, 9).Value = "1) 안전장비"
$ws.Cells.Item(3, 10).Value = "기타"
$ws.Cells.Item(3, 11).Value = "안전 1차 - 안전모내피"
$ws.Cells.Item(3, 12).Value = "KS규격-2"
$ws.Cells.Item(3, 13).Value = 5
$ws.Cells.Item(3, 14).Value = 2500
$ws.Cells.Item(3, 15).Value = 13750

# ---------------------------------------------------------------------
# Sheets "갑지" / "을지": clear the leftover empty 비고 cells (I2 / I3)
# ---------------------------------------------------------------------
foreach ($sheetName in @("갑지", "을지")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Cells.Item(2, 9).ClearContents()
    $sheet.Cells.Item(3, 9).ClearContents()
}
